# Updated cryptos list on Thu Jul 13 12:50:37 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.589.20'
$ws.Range("E2").Value = '  -0.71%  '

$ws.Range("D3").Value = '1.883.80'
$ws.Range("E3").Value = '  -0.40%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.39%  '

$ws.Range("E6").Value = '  -0.03%  '

$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2898'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.26%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06536'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.05%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.43'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.37%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7696'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.11%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '100.46'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.56%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07825'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.30%  '

$ws.Range("D14").Value = '1.884.75'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.244'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.35%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '285.12'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.06%  '

$ws.Range("D17").Value = '30.555.47'
$ws.Range("E17").Value = '  -2.31%  '

$ws.Range("E18").Value = '  -0.58%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007528'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.23%  '

$ws.Range("E20").Value = '  -0.02%  '

$ws.Range("D21").Value = '2.132.49'
$ws.Range("E21").Value = '  -0.41%  '

$ws.Range("E22").Value = '  -0.27%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.18%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.447'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.66%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.178'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.44%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.52'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.73%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.10'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.24%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.904'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.91%  '

$ws.Range("B29").Value = 'Stellar'
$ws.Range("C29").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.09698'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.58%  '

$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.329'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.01%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.501'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.65%  '

$ws.Range("E32").Value = '  -1.63%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.184'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.11%  '

$ws.Range("E34").Value = '  -0.73%  '

$ws.Range("E35").Value = '  -0.22%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6994'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.18%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.740'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.39%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01909'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.43%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.881'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.19%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '75.91'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.95%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.303'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.64%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.983'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.19%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4260'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.93%  '

$ws.Range("E44").Value = '  -0.10%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8319'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.82%  '

$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.45'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.54%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.899'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.09%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.029'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.42%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '35.16'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.33%  '

$ws.Range("E50").Value = '  +0.14%  '

$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3966'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.21%  '
